$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 73

# Column A holds the date as plain text in this sheet (e.g. "10/28/2025" in
# the row above), so enter it with a leading apostrophe to force text
# storage instead of having it auto-converted into a date serial number,
# then restore the default "Normal" style so no extra formatting is added.
$cellA = $ws.Cells.Item($row, 1)
$cellA.Value = "'10/29/2025"
$cellA.Style = "Normal"

# Column B: the day's numeric profit value.
$ws.Cells.Item($row, 2).Value = 11983.44
